# Commit: fix sending sms bugs
# Adds a new worksheet "2022_3" (after "2022_2") containing SMS/OTP send-otp
# log entries for Mon Mar 07 2022, mirroring the "home"/"2022_2" log-sheet layout.

$wb = $excel.ActiveWorkbook

# Insert the new sheet right after the last existing sheet (end of the tab strip).
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$ws.Name = "2022_3"

# Header row
$ws.Cells.Item(1, 1).Value = 'date'
$ws.Cells.Item(1, 2).Value = 'time'
$ws.Cells.Item(1, 3).Value = 'phoneNumber'
$ws.Cells.Item(1, 4).Value = 'model'
$ws.Cells.Item(1, 5).Value = 'path'
$ws.Cells.Item(1, 6).Value = 'action'
$ws.Cells.Item(1, 7).Value = 'status'
$ws.Cells.Item(1, 8).Value = 'description'
$ws.Cells.Item(1, 9).Value = 'failureReason'
$ws.Cells.Item(1, 10).Value = 'userId'
$ws.Cells.Item(1, 11).Value = 'modelId'

# Data rows (r2:r11)
# Row 2
$ws.Cells.Item(2, 1).Value = 'Mon Mar 07 2022'
$ws.Cells.Item(2, 2).Value = '07:20:39 GMT+0000 (Greenwich Mean Time)'
$ws.Cells.Item(2, 3).Value = "'" + '+22892942601'
$ws.Cells.Item(2, 4).Value = 'User'
$ws.Cells.Item(2, 5).Value = '/api/auth/send-otp'
$ws.Cells.Item(2, 6).Value = 'request'
$ws.Cells.Item(2, 7).Value = 'failed'
$ws.Cells.Item(2, 8).Value = '+22892942601 request to receive otp'
$ws.Cells.Item(2, 9).Value = 'phoneNumber.chartAt is not a function'
$ws.Cells.Item(2, 10).Value = '+22892942601 request to receive otp'
$ws.Cells.Item(2, 11).Value = 'phoneNumber.chartAt is not a function'
# Row 3
$ws.Cells.Item(3, 1).Value = 'Mon Mar 07 2022'
$ws.Cells.Item(3, 2).Value = '07:25:45 GMT+0000 (Greenwich Mean Time)'
$ws.Cells.Item(3, 3).Value = "'" + '+22892942601'
$ws.Cells.Item(3, 4).Value = 'User'
$ws.Cells.Item(3, 5).Value = '/api/auth/send-otp'
$ws.Cells.Item(3, 6).Value = 'request'
$ws.Cells.Item(3, 7).Value = 'failed'
$ws.Cells.Item(3, 8).Value = '+22892942601 request to receive otp'
$ws.Cells.Item(3, 9).Value = 'phoneNumber.chartAt is not a function'
# Row 4
$ws.Cells.Item(4, 1).Value = 'Mon Mar 07 2022'
$ws.Cells.Item(4, 2).Value = '07:27:09 GMT+0000 (Greenwich Mean Time)'
$ws.Cells.Item(4, 3).Value = "'" + '+22892942601'
$ws.Cells.Item(4, 4).Value = 'User'
$ws.Cells.Item(4, 5).Value = '/api/auth/send-otp'
$ws.Cells.Item(4, 6).Value = 'request'
$ws.Cells.Item(4, 7).Value = 'failed'
$ws.Cells.Item(4, 8).Value = '+22892942601 request to receive otp'
$ws.Cells.Item(4, 9).Value = 'phoneNumber.chartAt is not a function'
# Row 5
$ws.Cells.Item(5, 1).Value = 'Mon Mar 07 2022'
$ws.Cells.Item(5, 2).Value = '07:28:23 GMT+0000 (Greenwich Mean Time)'
$ws.Cells.Item(5, 3).Value = "'" + '+22892942601'
$ws.Cells.Item(5, 4).Value = 'User'
$ws.Cells.Item(5, 5).Value = '/api/auth/send-otp'
$ws.Cells.Item(5, 6).Value = 'request'
$ws.Cells.Item(5, 7).Value = 'failed'
$ws.Cells.Item(5, 8).Value = '+22892942601 request to receive otp'
$ws.Cells.Item(5, 9).Value = 'Cannot read properties of undefined (reading ''apiUrl'')'
# Row 6
$ws.Cells.Item(6, 1).Value = 'Mon Mar 07 2022'
$ws.Cells.Item(6, 2).Value = '07:30:00 GMT+0000 (Greenwich Mean Time)'
$ws.Cells.Item(6, 3).Value = "'" + '+22892942601'
$ws.Cells.Item(6, 4).Value = 'User'
$ws.Cells.Item(6, 5).Value = '/api/auth/send-otp'
$ws.Cells.Item(6, 6).Value = 'request'
$ws.Cells.Item(6, 7).Value = 'succeeded'
$ws.Cells.Item(6, 8).Value = '+22892942601 request to receive otp'
# Row 7
$ws.Cells.Item(7, 1).Value = 'Mon Mar 07 2022'
$ws.Cells.Item(7, 2).Value = '07:31:44 GMT+0000 (Greenwich Mean Time)'
$ws.Cells.Item(7, 3).Value = "'" + '22892942601'
$ws.Cells.Item(7, 4).Value = 'User'
$ws.Cells.Item(7, 5).Value = '/api/auth/send-otp'
$ws.Cells.Item(7, 6).Value = 'request'
$ws.Cells.Item(7, 7).Value = 'failed'
$ws.Cells.Item(7, 8).Value = '22892942601 request to receive otp'
$ws.Cells.Item(7, 9).Value = 'error.invalid'
# Row 8
$ws.Cells.Item(8, 1).Value = 'Mon Mar 07 2022'
$ws.Cells.Item(8, 2).Value = '07:34:17 GMT+0000 (Greenwich Mean Time)'
$ws.Cells.Item(8, 3).Value = "'" + '22892942601'
$ws.Cells.Item(8, 4).Value = 'User'
$ws.Cells.Item(8, 5).Value = '/api/auth/send-otp'
$ws.Cells.Item(8, 6).Value = 'request'
$ws.Cells.Item(8, 7).Value = 'succeeded'
$ws.Cells.Item(8, 8).Value = '22892942601 request to receive otp'
# Row 9
$ws.Cells.Item(9, 1).Value = 'Mon Mar 07 2022'
$ws.Cells.Item(9, 2).Value = '07:38:55 GMT+0000 (Greenwich Mean Time)'
$ws.Cells.Item(9, 3).Value = "'" + '22892942601'
$ws.Cells.Item(9, 4).Value = 'User'
$ws.Cells.Item(9, 5).Value = '/api/auth/send-otp'
$ws.Cells.Item(9, 6).Value = 'request'
$ws.Cells.Item(9, 7).Value = 'succeeded'
$ws.Cells.Item(9, 8).Value = '22892942601 request to receive otp'
# Row 10
$ws.Cells.Item(10, 1).Value = 'Mon Mar 07 2022'
$ws.Cells.Item(10, 2).Value = '07:41:48 GMT+0000 (Greenwich Mean Time)'
$ws.Cells.Item(10, 3).Value = "'" + '22892942601'
$ws.Cells.Item(10, 4).Value = 'User'
$ws.Cells.Item(10, 5).Value = '/api/auth/send-otp'
$ws.Cells.Item(10, 6).Value = 'request'
$ws.Cells.Item(10, 7).Value = 'succeeded'
$ws.Cells.Item(10, 8).Value = '22892942601 request to receive otp'
# Row 11
$ws.Cells.Item(11, 1).Value = 'Mon Mar 07 2022'
$ws.Cells.Item(11, 2).Value = '07:46:33 GMT+0000 (Greenwich Mean Time)'
$ws.Cells.Item(11, 3).Value = "'" + '22892942601'
$ws.Cells.Item(11, 4).Value = 'User'
$ws.Cells.Item(11, 5).Value = '/api/auth/send-otp'
$ws.Cells.Item(11, 6).Value = 'request'
$ws.Cells.Item(11, 7).Value = 'succeeded'
$ws.Cells.Item(11, 8).Value = '22892942601 request to receive otp'

# Restore the original active tab ("home") instead of leaving the newly
# inserted sheet selected.
$wb.Worksheets.Item(1).Select()

Write-Output ("Added sheet " + $ws.Name + " with " + $ws.UsedRange.Rows.Count + " rows")
